$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82-140 down to 83-141
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new record
$ws.Cells.Item(82,1).Value = 8
$ws.Cells.Item(82,2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(82,3).Value = 'Coquimbo'
$ws.Cells.Item(82,4).Value = 44977
$ws.Cells.Item(82,5).Value = 4
$ws.Cells.Item(82,6).Value = 'Fruta'
$ws.Cells.Item(82,7).Value = 100109
$ws.Cells.Item(82,8).Value = 'Uva'
$ws.Cells.Item(82,9).Value = 100109001
$ws.Cells.Item(82,10).Value = 'Uva'
$ws.Cells.Item(82,11).Value = 'Ralli Seedless'
$ws.Cells.Item(82,12).Value = 'Primera'
$ws.Cells.Item(82,13).Value = 240
$ws.Cells.Item(82,14).Value = 10000
$ws.Cells.Item(82,15).Value = 11000
$ws.Cells.Item(82,16).Value = 10500
$ws.Cells.Item(82,17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(82,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(82,19).Value = 583
$ws.Cells.Item(82,20).Value = 18
